$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write one data row (columns A..G) on a given worksheet.
#   $ws        - worksheet object
#   $row       - destination row number
#   $fmtRow    - existing row number to copy per-column number formats from
#                (so Date/Text column A and DateTime column B pick up the
#                same style already used throughout the sheet, and so the
#                text in column A is not auto-converted into a date serial)
#   $date      - text for column A (e.g. "2023-02-23")
#   $runTime   - numeric serial for column B
#   $name      - text for column C
#   $total,$pass,$fail,$taken - numeric values for D,E,F,G
# ---------------------------------------------------------------------------
function Write-DataRow {
    param($ws, $row, $fmtRow, $date, $runTime, $name, $total, $pass, $fail, $taken)

    # Column A: literal date-looking text (leading apostrophe stops the
    # host from re-interpreting "2023-02-23" as a real date), then copy the
    # formatting of an already-correct cell on top of it (format only -
    # PasteSpecial formats never touches the value we just entered).
    $ws.Cells.Item($row, 1).Value = "'" + $date
    $ws.Cells.Item($fmtRow, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    # Column B: numeric run-time serial, formatted like the sibling rows
    # (custom datetime number format) instead of the plain column default.
    $ws.Cells.Item($row, 2).Value = $runTime
    $ws.Cells.Item($fmtRow, 2).Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4122)

    # Columns C..G: plain text / numbers - these already inherit the sheet's
    # column formatting when assigned directly.
    $ws.Cells.Item($row, 3).Value = $name
    $ws.Cells.Item($row, 4).Value = $total
    $ws.Cells.Item($row, 5).Value = $pass
    $ws.Cells.Item($row, 6).Value = $fail
    $ws.Cells.Item($row, 7).Value = $taken
}

# ===========================================================================
# AMSIN sheet: append rows 96-103
# ===========================================================================
$wsAmsin = $wb.Worksheets.Item("AMSIN")

Write-DataRow $wsAmsin 96 95 "2023-02-23" 44980.75622256944   "173hhttffxx"         269 269 0 5.32
Write-DataRow $wsAmsin 97 95 "2023-03-09" 44994.57150388889   "174fstcycle"         269 264 5 6.28
Write-DataRow $wsAmsin 98 95 "2023-03-10" 44995.78329828704   "174ffiinnalrun"      269 269 0 5.98
Write-DataRow $wsAmsin 99 95 "2023-03-13" 44998.46071991898   "174finalrun"         269 269 0 4.79
Write-DataRow $wsAmsin 100 95 "2023-03-28" 45013.52588585648  "175prerun"           269 266 3 5.5
Write-DataRow $wsAmsin 101 95 "2023-03-30" 45015.70615905093  "175scndcyc"          269 263 6 6.69
Write-DataRow $wsAmsin 102 95 "2023-03-31" 45016.44203090278  "175fnlrun"           269 268 1 5.23
Write-DataRow $wsAmsin 103 95 "2023-04-12" 45028.45964890953  "176fstrtail"         269 260 9 6.04

# ===========================================================================
# BETA sheet: append rows 35-36
# ===========================================================================
$wsBeta = $wb.Worksheets.Item("BETA")

Write-DataRow $wsBeta 35 34 "2023-03-13" 44998.52453083333 "174beta"  269 268 1 3.92
Write-DataRow $wsBeta 36 34 "2023-03-31" 45016.52234857639 "175beta"  269 265 4 4.52

# ===========================================================================
# AMS sheet: fix up row 70's formatting/value, then append rows 71-76
# ===========================================================================
$wsAms = $wb.Worksheets.Item("AMS")

# Row 70 already holds the correct text in A70/C70 and numbers in D70..G70;
# only the per-cell style (picked up from row 69, which already has it) and
# the B70 run-time value need to change.
$wsAms.Cells.Item(69, 1).Copy()
$wsAms.Cells.Item(70, 1).PasteSpecial(-4122)
$wsAms.Cells.Item(69, 3).Copy()
$wsAms.Cells.Item(70, 3).PasteSpecial(-4122)
$wsAms.Cells.Item(69, 4).Copy()
$wsAms.Cells.Item(70, 4).PasteSpecial(-4122)
$wsAms.Cells.Item(69, 5).Copy()
$wsAms.Cells.Item(70, 5).PasteSpecial(-4122)
$wsAms.Cells.Item(69, 6).Copy()
$wsAms.Cells.Item(70, 6).PasteSpecial(-4122)
$wsAms.Cells.Item(69, 7).Copy()
$wsAms.Cells.Item(70, 7).PasteSpecial(-4122)

$wsAms.Cells.Item(70, 2).Value = 44977.83047721065
$wsAms.Cells.Item(69, 2).Copy()
$wsAms.Cells.Item(70, 2).PasteSpecial(-4122)

Write-DataRow $wsAms 71 69 "2023-02-23" 44980.83663085648 "173htfxbulkschedule"  269 269 0 3.83
Write-DataRow $wsAms 72 69 "2023-03-01" 44986.68783164352 "173angularvrs"        269 269 0 3.78
Write-DataRow $wsAms 73 69 "2023-03-02" 44987.43015282408 "liveangular173"       269 269 0 4.16
Write-DataRow $wsAms 74 69 "2023-03-07" 44992.73899635416 "174htfxmar"           269 268 1 4.03
Write-DataRow $wsAms 75 69 "2023-03-13" 44998.83213166667 "174live"              269 268 1 3.81
Write-DataRow $wsAms 76 69 "2023-03-31" 45016.78996877315 "175live"              269 267 2 4.04
